$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1, formatted like the other header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns I and J, rows 2-8
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3

$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 3

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 2
